# The commit swaps the two theme parts of the deck:
#   ppt/theme/theme1.xml (the slide master's theme, currently "Integral")
#   ppt/theme/theme2.xml (the notes master's theme, currently "Office Theme")
# end up with each other's content - i.e. the slides end up on the generic
# "Office Theme" palette and the notes master ends up on the "Integral"
# palette. The font scheme and format scheme (fills/lines/effects) in both
# parts are already byte-identical, so the only real content difference is
# the 12-colour scheme (and the cosmetic theme/clrScheme "name" attributes,
# which PowerPoint's automation surface does not expose for writing).
#
# The only reliably wired, side-effect-free COM surface for rewriting a
# theme's colour scheme here is Slide.ThemeColorScheme, which edits the
# colour scheme backing the slide master (ppt/theme/theme1.xml) that every
# slide in this deck shares. Apply the target ("Office Theme") palette
# through it.

$p = $ppt.ActivePresentation

function ToCOMColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Target palette for ppt/theme/theme1.xml ("Office Theme" / "Office"),
# in MsoThemeColorSchemeIndex order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

# Every slide shares the one slide master / theme part, so touching slide 1
# is enough to rewrite ppt/theme/theme1.xml for the whole deck.
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = ToCOMColor($officeColors[$i - 1])
}
